$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new activity rows below the existing scoreboard data (last existing row is 133).
# Columns: A Participant, B Date, C Workout Type, D Total Duration, E Total Distance,
#          F Total Elevation, G Zone1, H Zone2, I Zone3, J Zone4, K Zone5, L Workout Level, M Week

$rows = @(
    @("Matt", 45473, "Run",     40, 3.52, 217, 3,  23, 8,  2, 0, "Agile Antelope",   3),
    @("Matt", 45473, "Walk",    5,  0.2,  13,  5,  0,  0,  0, 0, "Agile Antelope",   3),
    @("Eric", 45473, "Workout", 32, 0,    0,   32, 0,  0,  0, 0, "Sauntering Hippo", 3)
)

$startRow = 134
$lastExistingRow = 133

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting (incl. the date number format) from the last existing data row first,
    # so the new rows reuse the workbook's existing cell styles instead of minting new ones.
    $ws.Range($ws.Cells.Item($lastExistingRow, 1), $ws.Cells.Item($lastExistingRow, 13)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 13)).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
    $ws.Cells.Item($r, 13).Value = $data[12]
}

$excel.CutCopyMode = $false

# Move the selection to reflect where a user would land after entering the new rows.
$ws.Range("A137").Select()
